$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $helper = $ws.Range("ZZ1")
    $helper.NumberFormat = "@"
    $helper.Value = $value
    $helper.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $helper.Clear()
}

Set-TextValue 'D2' '44.509.82'
Set-TextValue 'D3' '2.239.97'
Set-TextValue 'E3' '  -0.22%  '
Set-TextValue 'E4' '  +0.34%  '
Set-TextValue 'D5' '304.77'
Set-TextValue 'E5' '  -0.75%  '
Set-TextValue 'D6' '94.66'
Set-TextValue 'E6' '  -1.11%  '
Set-TextValue 'D7' '0.569'
Set-TextValue 'E7' '  -0.81%  '
Set-TextValue 'E8' '  +0.16%  '
Set-TextValue 'E9' '  -2.37%  '
Set-TextValue 'D10' '34.74'
Set-TextValue 'E10' '  -1.07%  '
Set-TextValue 'D11' '0.0799'
Set-TextValue 'E11' '  -2.00%  '
Set-TextValue 'E12' '  -1.61%  '
Set-TextValue 'E13' '  -0.22%  '
Set-TextValue 'D14' '2.581.73'
Set-TextValue 'E14' '  -0.23%  '
Set-TextValue 'D15' '2.235.10'
Set-TextValue 'E15' '  -4.16%  '
Set-TextValue 'D16' '0.828'
Set-TextValue 'E16' '  -0.71%  '
Set-TextValue 'D17' '13.47'
Set-TextValue 'E17' '  -1.11%  '
Set-TextValue 'D18' '44.314.86'
Set-TextValue 'E18' '  +0.58%  '
Set-TextValue 'E19' '  -3.39%  '
Set-TextValue 'D20' '6.17'
Set-TextValue 'E20' '  -3.78%  '
Set-TextValue 'D21' '11.73'
Set-TextValue 'E21' '  -3.71%  '
Set-TextValue 'D22' '64.98'
Set-TextValue 'E22' '  -0.87%  '
Set-TextValue 'D23' '236.96'
Set-TextValue 'E23' '  -0.02%  '
Set-TextValue 'E24' '  -1.39%  '
Set-TextValue 'E25' '  -2.39%  '
Set-TextValue 'E26' '  -0.07%  '
Set-TextValue 'E27' '  +5.04%  '
Set-TextValue 'D28' '9.67'
Set-TextValue 'E28' '  -3.25%  '
Set-TextValue 'D29' '37.10'
Set-TextValue 'E29' '  -2.68%  '
Set-TextValue 'D30' '19.82'
Set-TextValue 'E30' '  -1.73%  '
Set-TextValue 'D31' '5.81'
Set-TextValue 'E31' '  -2.27%  '
Set-TextValue 'D32' '149.99'
Set-TextValue 'E32' '  -1.90%  '
Set-TextValue 'B33' 'Hedera'
Set-TextValue 'C33' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D33' '0.0780'
Set-TextValue 'E33' '  -2.55%  '
Set-TextValue 'B34' 'WEMIXToken'
Set-TextValue 'C34' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D34' '2.62'
Set-TextValue 'E34' '  +0.39%  '
Set-TextValue 'E35' '  -2.09%  '
Set-TextValue 'D36' '1.88'
Set-TextValue 'E36' '  +6.98%  '
Set-TextValue 'B37' 'Stellar'
Set-TextValue 'C37' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D37' '0.118'
Set-TextValue 'E37' '  -2.11%  '
Set-TextValue 'B38' 'Kaspa'
Set-TextValue 'C38' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D38' '0.107'
Set-TextValue 'E38' '  -0.93%  '
Set-TextValue 'D39' '14.83'
Set-TextValue 'E39' '  +1.61%  '
Set-TextValue 'E40' '  -5.46%  '
Set-TextValue 'E41' '  -2.41%  '
Set-TextValue 'D42' '0.0296'
Set-TextValue 'E42' '  -0.96%  '
Set-TextValue 'E43' '  +0.17%  '
Set-TextValue 'D44' '1.820.61'
Set-TextValue 'E44' '  +3.60%  '
Set-TextValue 'D45' '1.75'
Set-TextValue 'E45' '  +10.40%  '
Set-TextValue 'D46' '79.41'
Set-TextValue 'E46' '  -4.56%  '
Set-TextValue 'D47' '0.187'
Set-TextValue 'E47' '  -2.72%  '
Set-TextValue 'D48' '98.01'
Set-TextValue 'E48' '  -2.36%  '
Set-TextValue 'D49' '4.83'
Set-TextValue 'E49' '  -2.61%  '
Set-TextValue 'D50' '68.44'
Set-TextValue 'E50' '  +0.77%  '
Set-TextValue 'E51' '  -2.39%  '

